# Problem with material balances.
# Insert a new "Type_2" column into BaseFlows (sheet2), shifting the old
# "feed" column from D to E, and populate the new column with the
# material-balance stage for each flow. Also re-point the active sheet/
# selection/zoom to BaseFlows, matching the author's final view state.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("BaseFlows")

# --- Insert a new column before D (old D -> E, keeps old formatting) ---
$ws.Columns("D:D").Insert()

# The inserted column initially inherits formatting from column C; clear
# that back to the workbook default for the data rows before writing the
# new values (the header cell D1 keeps the inherited bold header style,
# same as the rest of row 1).
$ws.Range("D2:D40").ClearFormats()

# --- Header ---
$ws.Range("D1").Value = "Type_2"

# --- Data rows: the processing stage each flow's output corresponds to ---
$ws.Range("D2").Value = "Iron ore"
$ws.Range("D3").Value = "Reduced iron"
$ws.Range("D4").Value = "Reduced iron"
$ws.Range("D5").Value = "Liquid steel"
$ws.Range("D6").Value = "Liquid steel"
$ws.Range("D7").Value = "Liquid steel"
$ws.Range("D8").Value = "Castings"
$ws.Range("D9").Value = "Castings"
$ws.Range("D10").Value = "Castings"
$ws.Range("D11").Value = "Castings"
$ws.Range("D12").Value = "Finished steel"
$ws.Range("D13").Value = "Finished steel"
$ws.Range("D14").Value = "Finished steel"
$ws.Range("D15").Value = "Finished steel"
$ws.Range("D16").Value = "Finished steel"
$ws.Range("D17").Value = "Finished steel"
$ws.Range("D18").Value = "Finished steel"
$ws.Range("D19").Value = "Finished steel"
$ws.Range("D20").Value = "Finished steel"
$ws.Range("D21").Value = "Finished steel"
$ws.Range("D22").Value = "Finished steel"
$ws.Range("D23").Value = "Finished steel"
$ws.Range("D24").Value = "Finished steel"
$ws.Range("D25").Value = "Finished steel"
$ws.Range("D26").Value = "Finished steel"
$ws.Range("D27").Value = "Finished steel"
$ws.Range("D28").Value = "Final product"
$ws.Range("D29").Value = "Final product"
$ws.Range("D30").Value = "Final product"
$ws.Range("D31").Value = "Final product"
$ws.Range("D32").Value = "Final product"
$ws.Range("D33").Value = "Final product"
$ws.Range("D34").Value = "Final product"
$ws.Range("D35").Value = "Final product"
$ws.Range("D36").Value = "Final product"
$ws.Range("D37").Value = "Final product"
$ws.Range("D38").Value = "Scrap"
$ws.Range("D39").Value = "Scrap"
$ws.Range("D40").Value = "Scrap"

# --- Column width for the new column ---
$ws.Columns("D:D").ColumnWidth = 11.436197916666666

# --- Make BaseFlows the active sheet/tab, matching the saved view state ---
$ws.Activate()
$excel.ActiveWindow.Zoom = 85
$ws.Range("D39").Select()
